$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the numeric-looking answer columns (B:D for rows 2-3) as Text so
# that values such as "4.259" or "26" are stored as shared-string text
# (matching the original workbook's convention of keeping every cell as a
# string) instead of being auto-converted into numeric cells.
$ws.Range("B2:D3").NumberFormat = "@"

# Row 2: new question / tokens / time / answer
$ws.Range("A2").Value = "What is the average number of runs Mumbai have made in wins in overs 1-6 in the 2024 IPL?"
$ws.Range("B2").Value = "4.259"
$ws.Range("C2").Value = "26"
$ws.Range("D2").Value = "1.63"

# Row 3: repeated question / tokens / time / answer
$ws.Range("A3").Value = "What is the average number of runs Mumbai have made in wins in overs 1-6 in the 2024 IPL? - > this question repeatedIeB"
$ws.Range("B3").Value = "2.376"
$ws.Range("C3").Value = "5"
$ws.Range("D3").Value = "1.78"

# Remove the temporary Text number format again so the cells end up with the
# same (default) style they started with.
$ws.Range("B2:D3").ClearFormats()
